# The updated workbook drops the two "% Change" columns (originally column M,
# "% Change from 2013-2014 to 2014-2015", and column X, "% Change from
# 2023-2024 to 2024-2025") from the single data sheet. Deleting column X
# first means column M's letter reference is still valid afterwards, so the
# two deletes don't need any re-computed offsets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("X").Delete()
$ws.Columns("M").Delete()
